$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph ("Bivariate Pearson correlations revealed significant ...") --
# rewrite to add the "95% CI [...] did not cross zero" sentence, move the
# 95% CI figures into the r-value sentences, and drop the separate p-value
# clauses (merging their surrounding prose into the neighbouring plain runs).
# ---------------------------------------------------------------------------

function Get-FoundRange($searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $searchText"
    }
    return $rng
}

# 1) Opening sentence: add the new CI sentence before "Negative affect ..."
#    and change the trailing "(" into ", " (the r value no longer needs a
#    parenthesis since the stats now read "r = ..., 95% CI [...]").
$rng = Get-FoundRange("Bivariate Pearson correlations revealed significant associations among the continuous variables. Negative affect was significantly negatively correlated with conscientiousness (")
$rng.Text = "Bivariate Pearson correlations revealed significant associations among the continuous variables, all supported by 95% confidence intervals that did not cross zero, indicating that the effects were likely meaningful and not due to chance. Negative affect was significantly negatively correlated with conscientiousness, "

# 2) First r value: " = -.37, " gains its 95% CI and the "suggesting..." clause.
$rng = Get-FoundRange(" = –.37, ")
$rng.Text = " = –.37, 95% CI [–.50, –.21], suggesting that individuals who rated themselves as more conscientious tended to report lower levels of negative affect. "

# 3) Remove the first "p" run (" < .01) ... SPP was positively associated with
#    negative affect (") and replace the trailing prose with the shorter,
#    comma-joined phrasing that now leads into the second r value.
$rng = Get-FoundRange(" < .01), indicating that individuals with higher conscientiousness tended to report lower levels of negative affect. In contrast, SPP was positively associated with negative affect (")
$blockStart = $rng.Start
$blockEnd = $rng.End
$pRange = $d.Range($blockStart - 1, $blockStart)
$pRange.Delete()
$blockRange = $d.Range($blockStart - 1, $blockEnd - 1)
$blockRange.Text = "SPP was positively associated with negative affect, "

# 4) Second r value: " = .39, " gains its 95% CI and the "indicating..." clause.
$rng = Get-FoundRange(" = .39, ")
$rng.Text = " = .39, 95% CI [.24, .52], indicating that participants who perceived stronger external expectations to be perfect also reported more negative affect. Finally, conscientiousness and SPP were modestly but significantly negatively correlated, "

# 5) Remove the second "p" run (" < .01) ... correlated (") and replace the
#    trailing prose - nothing remains to insert here since the lead-in text
#    was already appended above, so this clause collapses to nothing.
$rng = Get-FoundRange(" < .01), suggesting that participants with stronger beliefs that others expect them to be perfect also reported more negative affect. Finally, conscientiousness and SPP were modestly but significantly negatively correlated (")
$blockStart = $rng.Start
$blockEnd = $rng.End
$pRange = $d.Range($blockStart - 1, $blockStart)
$pRange.Delete()
$blockRange = $d.Range($blockStart - 1, $blockEnd - 1)
$blockRange.Text = ""

# 6) Third r value: " = -.20, " gains its 95% CI.
$rng = Get-FoundRange(" = –.20, ")
$rng.Text = " = –.20, 95% CI [–.36, –.04], "

# 7) Remove the third "p" run (" < .05) ... traits.") - the remaining
#    "reflecting some conceptual..." prose is retained verbatim.
$rng = Get-FoundRange(" < .05), reflecting some conceptual overlap but also distinctiveness between these personality traits.")
$blockStart = $rng.Start
$blockEnd = $rng.End
$pRange = $d.Range($blockStart - 1, $blockStart)
$pRange.Delete()
$blockRange = $d.Range($blockStart - 1, $blockEnd - 1)
$blockRange.Text = "reflecting some conceptual overlap but also distinctiveness between these personality traits."

Write-Output "paragraph rewrite complete"
